$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.141481995582581
$ws.Range("B1").Value = 2.232837915420532
$ws.Range("C1").Value = 11.07173156738281
$ws.Range("D1").Value = 2.332935571670532
$ws.Range("E1").Value = 1.274168729782104
